# Add Arabic guidance row to the Skid sheet (row 2), pushing the
# existing data row down to row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the data row.
$ws.Rows.Item(2).Insert()

# Fill in the Arabic guidance text for the three columns.
$ws.Range("A2").Value = "أدخل رمز القطاع/المقطع"
$ws.Range("B2").Value = "أدخل رقم المسار (L1, L2, إلخ)"
$ws.Range("C2").Value = "أدخل قيمة معامل الاحتكاك/مقاومة الانزلاق"

# Build the guidance-row look (italic, small, grey font on a light-grey
# fill, right aligned / vertically centered / wrapped) on a scratch cell
# far away from the used range, then copy *just the formatting* onto the
# guidance row so it becomes a single new style instead of one new style
# per property assignment.
$helper = $ws.Range("ZZ9999")
$helper.Font.Italic = $true
$helper.Font.Size = 9
$helper.Font.Color = 6710886
$helper.Interior.Color = 15790320
$helper.HorizontalAlignment = -4152
$helper.VerticalAlignment = -4108
$helper.WrapText = $true

$helper.Copy()
$guidanceRange = $ws.Range("A2:C2")
$guidanceRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$helper.Clear()

$guidanceRange.RowHeight = 30

Write-Host "Guidance row inserted."
